$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts old D:K data to E:L), matching the
# new fiscal-year column added to the Income Statement / Balance Sheet /
# Cash Flow Statement tables.
$ws.Columns("D").Insert()

# Carry over the number formats/styles from the column that used to be D
# (now E) onto the freshly inserted column D so date/number formatting
# matches the rest of the table.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the new period's figures.
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 1606200
$ws.Range("D9").Value2 = 859900
$ws.Range("D10").Value2 = 746300
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 45000
$ws.Range("D15").Value2 = 185000
$ws.Range("D17").Value2 = 1371400
$ws.Range("D18").Value2 = 234800
$ws.Range("D20").Value2 = -400
$ws.Range("D21").Value2 = 419400
$ws.Range("D22").Value2 = 125700
$ws.Range("D23").Value2 = 108700
$ws.Range("D24").Value2 = 4900
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 103800
$ws.Range("D27").Value2 = 107900
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = "NA"
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = 400
$ws.Range("D33").Value2 = 107900
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 107900
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 52700
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 264900
$ws.Range("D44").Value2 = 0
$ws.Range("D45").Value2 = 111900
$ws.Range("D46").Value2 = 429500
$ws.Range("D47").Value2 = 0
$ws.Range("D48").Value2 = 652900
$ws.Range("D49").Value2 = 2616900
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 129400
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 3828700
$ws.Range("D57").Value2 = 68400
$ws.Range("D58").Value2 = 160000
$ws.Range("D59").Value2 = 174200
$ws.Range("D60").Value2 = 402600
$ws.Range("D61").Value2 = 2149600
$ws.Range("D62").Value2 = 131200
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 2725900
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = -871600
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 1102800
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 107900
$ws.Range("D83").Value2 = 185000
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 214300
$ws.Range("D91").Value2 = -82300
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -90400
$ws.Range("D96").Value2 = -203900
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = -117700
$ws.Range("D101").Value2 = -400
$ws.Range("D102").Value2 = 5800
"done"
